$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.922.06'
$ws.Range('E2').Value = '  +0.20%  '
$ws.Range('D3').Value = '2.525.40'
$ws.Range('E3').Value = '  -2.63%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '589.28'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.97%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '172.62'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +3.99%  '
$ws.Range('E7').Value = '  -0.05%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.527'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.03%  '
$ws.Range('D9').Value = '2.524.71'
$ws.Range('E9').Value = '  -2.37%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.137'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.82%  '
$ws.Range('E11').Value = '  +1.81%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '5.12'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -1.26%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.343'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -6.26%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '26.66'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.75%  '
$ws.Range('D15').Value = '3.009.73'
$ws.Range('E15').Value = '  -1.88%  '
$ws.Range('E16').Value = '  -1.79%  '
$ws.Range('D17').Value = '66.793.81'
$ws.Range('E17').Value = '  -0.52%  '
$ws.Range('D18').Value = '2.549.85'
$ws.Range('E18').Value = '  -2.45%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '8.01'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +3.15%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.27'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -2.95%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '353.43'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.08%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.17'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -2.00%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.59'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.67%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.98'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +4.47%  '
$ws.Range('E25').Value = '  +0.02%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '69.61'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.76%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.02'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -3.32%  '
$ws.Range('D28').Value = '2.679.23'
$ws.Range('E28').Value = '  -1.91%  '
$ws.Range('E29').Value = '  -0.13%  '
$ws.Range('D30').Value = '0.0₃0981'
$ws.Range('E30').Value = '  -0.73%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '531.47'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.72%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '8.12'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.33%  '
$ws.Range('E33').Value = '  +0.26%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.84'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.52%  '
$ws.Range('E35').Value = '  -1.55%  '
$ws.Range('E36').Value = '  +0.22%  '
$ws.Range('E37').Value = '  -0.94%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '157.38'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.23%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '18.59'
$ws.Range('D39').Style = 'Normal'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '18.43'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +1.03%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.354'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -2.01%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.78'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.38%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.10'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.37%  '
$ws.Range('E44').Value = '  -0.04%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.49'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +3.41%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '39.69'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.07%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '148.63'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.75%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.556'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -2.87%  '
$ws.Range('D49').Value = '0.0₆0277'
$ws.Range('E49').Value = '  -4.91%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '3.67'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.86%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.68'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.80%  '

Write-Host "Updated cryptos list"
